# feat: implement book addition and removal for admin interface
#
# The "Cărți favorite" sheet lists favourite books in rows 4-8 (row 2 is the
# title, row 3 the column headers). This edit:
#   - removes 3 books (rows that held "Adobe InDesign CC...", "Povesti in
#     romana si germana" and "Dictionar Oxford...")
#   - replaces the first book ("Drumul spre Biserica") with a newly added
#     book ("In apararea pietelor")
#   - removes the now-unused blank spacer row 1 and the decorative column A
#     formatting that ran down the left edge of the table

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove 3 books (rows 6, 7 and 8) -------------------------------------
$ws.Rows("6:8").Delete()

# --- Replace the "removed" book in row 4 with the newly added book -------
$ws.Range("B4").Value = "In apararea pietelor"
$ws.Range("C4").Value = "James Lucian"
$ws.Range("D4").Value = "Curtea Veche"
$ws.Range("E4").Value = "2012"
$ws.Range("F4").Value = "416"
$ws.Range("G4").Value = "Economie"
$ws.Range("H4").Value = "9786065883239"
$ws.Range("I4").Value = "185"

# --- Drop the decorative, empty column-A cells in the table rows ---------
$ws.Range("A2:A5").Clear()

# --- Drop the blank spacer row at the very top of the sheet --------------
$ws.Range("A1:I1").Clear()
